$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.608123302459717
$ws.Range("B1").Value = 3.858084201812744
$ws.Range("C1").Value = 7.523407459259033
$ws.Range("D1").Value = 7.679680824279785
$ws.Range("E1").Value = 6.009304046630859
